$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF column (F) values for several rows as part of
# re-pulling data / recomputing the mean
$ws.Range("F4").Value = -5
$ws.Range("F5").Value = -4
$ws.Range("F6").Value = -7
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = -1
